$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated classification metric values (Precision, Recall, Accuracy, F1-Score)
# for rows 2-8 (columns C, D, E, F)

$ws.Range("C2").Value = 0.6475450443515565
$ws.Range("D2").Value = 0.5304488830553094
$ws.Range("E2").Value = 0.6562876471133352
$ws.Range("F2").Value = 0.5531300020617894

$ws.Range("C3").Value = 0.6788555333890164
$ws.Range("D3").Value = 0.6330308881811546
$ws.Range("E3").Value = 0.6844592716152349
$ws.Range("F3").Value = 0.6495117009622303

$ws.Range("C4").Value = 0.6978232450277249
$ws.Range("D4").Value = 0.6613670820266967
$ws.Range("E4").Value = 0.7167083680845149
$ws.Range("F4").Value = 0.6607971017287977

$ws.Range("C5").Value = 0.8060669117162877
$ws.Range("D5").Value = 0.7715082882922708
$ws.Range("E5").Value = 0.7843573348160504
$ws.Range("F5").Value = 0.7794834848664745

$ws.Range("C6").Value = 0.8219834940351166
$ws.Range("D6").Value = 0.8018011998335816
$ws.Range("E6").Value = 0.8102122138819386
$ws.Range("F6").Value = 0.8106492623521262

$ws.Range("C7").Value = 0.7957164589409468
$ws.Range("D7").Value = 0.7804162539177379
$ws.Range("E7").Value = 0.7969604299879529
$ws.Range("F7").Value = 0.7872969362904352

$ws.Range("C8").Value = 0.7879161065596556
$ws.Range("D8").Value = 0.7512293710654735
$ws.Range("E8").Value = 0.7833379668242053
$ws.Range("F8").Value = 0.7670381257688823
